$d = $word.ActiveDocument

# The paragraph ends with two "error blocks" right after "::Invalid)}":
#   1) "    " + "<---" + "Expression "a" is invalid: invalid type literal notExisting::Invalid"   -> remove entirely
#   2) "    " + "<---" + "invalid type literal notExisting::Invalid"                                -> keep (unchanged)
# We must delete the whole first block (3 runs) while leaving the second block's
# runs ("<---" and "invalid type literal...") as two separate runs (not merged).

$target = "    <---Expression `"a`" is invalid: invalid type literal notExisting::Invalid"

$findRange = $d.Content
$found = $findRange.Find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $findRange.Delete()
}

# Deleting that span merges the remaining trailing runs ("<---" and
# "invalid type literal notExisting::Invalid") into a single run because they
# share identical formatting. Force them back apart into two distinct runs by
# briefly nudging the "<---" run's color to a different value (splitting the
# run) and then restoring the original color.
$arrowRange = $d.Content
$foundArrow = $arrowRange.Find.Execute("<---invalid type literal notExisting::Invalid", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundArrow) {
    $splitRange = $d.Range($arrowRange.Start, $arrowRange.Start + 4)
    $splitRange.Font.Color = 16711935
    $splitRange.Font.Color = 255
}
